$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A7").Value = "AEA / ASSA meetings"
$ws.Range("B7").Value = "changes every year"
$ws.Range("C7").Value = 43191
$ws.Range("C7").NumberFormat = "d-mmm"
$ws.Range("D7").Value = "January"

$ws.Range("A7").HorizontalAlignment = -4131
$ws.Range("B7").HorizontalAlignment = -4108
$ws.Range("C7").HorizontalAlignment = -4108
$ws.Range("D7").HorizontalAlignment = -4108

$ws.Range("A7:D7").Borders.LineStyle = 1

$ws.Range("C8").Select()
